# "Added 1.1.0 of term"
# Bumps the published Version of this FHIR term artifact and refreshes the
# Date stamp on the "Metadata" property sheet:
#   Version : 1.0.0 -> 1.1.0
#   Date    : 2023-06-07T11:52:14+02:00 -> 2023-07-10T23:08:03+02:00

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$labelColumn = $ws.Columns.Item(1)

$versionLabel = $labelColumn.Find("Version")
$versionLabel.Offset(0, 1).Value = "1.1.0"

$dateLabel = $labelColumn.Find("Date")
$dateLabel.Offset(0, 1).Value = "2023-07-10T23:08:03+02:00"
